# Remove three training-set rows that were pruned from the data set:
#   - Customer Communication / Notification Management
#   - Customer Communication / Email Communication
#   - Money Movement - Inbound / Principal Payment ("SOFR Term repayment" row)
#
# In the original workbook these are rows 66, 67 and 86 (1-indexed, row 1 is
# the header). Delete bottom-up so earlier row numbers stay valid.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(86).Delete()
$ws.Rows.Item(67).Delete()
$ws.Rows.Item(66).Delete()
